$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition / resize the logo picture (moved from top-right to top-left, slightly smaller) ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 9.708897637795275
$shp.Top = 6.00007874015748
$shp.Width = 80.57062992125984
$shp.Height = 60.27267716535433

# --- Move the title + table content further down the sheet ---
# old layout: D1 = title ("Bordereau d'envoi", bold/size13), row 3 = table header labels
# new layout: B7 = title, F7 = "le", row 11 = table header labels

# Stash the title cell's existing format (bold, size 13) in a scratch cell far away so it
# survives the row surgery below, then we can re-apply it without creating a duplicate style.
$ws.Range("D1").Copy()
$ws.Range("A1000").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Clear the old title cell and push the remaining rows down: deleting the 2 now-empty rows
# then inserting 10 blank ones lands the old row-3 table header exactly on row 11.
$ws.Range("D1").Clear()
$ws.Rows("1:2").Delete()
$ws.Rows("1:10").Insert()

# Re-apply the stashed title format to the new title cell, then clean up the scratch cell.
$ws.Range("A1008").Copy()
$ws.Range("B7").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A1008").Clear()

$ws.Range("B7").Value = "Bordereau d'envoi"
$ws.Rows("7").RowHeight = 17.4

# New "le" label alongside the title
$ws.Range("F7").Value = "le"

# --- Sheet view changes ---
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("B11").Select()
